$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C46").Value = 279.14
$ws.Range("D46").Value = 185
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Formula = "=(H46+I46)-(C46+D46+E46+F46+G46)"

$ws.Range("C47").Value = 450
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 17.7
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Formula = "=(H47+I47)-(C47+D47+E47+F47+G47)"

$ws.Range("J47").Select()
